$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet has a repeating 4-column block (Alain/Henri/Tony/Dulcinee header,
# OUI/NON answers below) tiled across many columns, followed by two trailing
# "summary" columns (email address, then an empty numeric column).
#
# We extend the tiling by 5 more repetitions (20 columns), inserting them
# right before the two trailing summary columns so the summary columns shift
# right and keep their own formatting/content.

$lastRow = 9
$blockWidth = 4
$repeats = 5
$newCols = $blockWidth * $repeats

# Column letters: the first tiled block (already styled, s="3") lives at
# columns I:L (9:12). The trailing summary columns are APA (1093) and
# APB (1094) before the edit.
$sourceStartCol = 9
$insertBeforeCol = 1093

# Insert blank columns immediately before the first summary column; Excel
# shifts the summary columns to the right and the newly inserted columns
# inherit the formatting of the column to their left (s="3"), matching the
# rest of the tiled block.
$insertStartLetter = $ws.Cells.Item(1, $insertBeforeCol).Address($false, $false)
$insertEndLetter = $ws.Cells.Item(1, $insertBeforeCol + $newCols - 1).Address($false, $false)
$ws.Range($insertStartLetter + ":" + $insertEndLetter).EntireColumn.Insert()

# Fill the newly inserted columns by tiling the existing 4-column pattern.
for ($row = 1; $row -le $lastRow; $row++) {
    for ($rep = 0; $rep -lt $repeats; $rep++) {
        for ($offset = 0; $offset -lt $blockWidth; $offset++) {
            $srcCol = $sourceStartCol + $offset
            $dstCol = $insertBeforeCol + ($rep * $blockWidth) + $offset
            $srcValue = $ws.Cells.Item($row, $srcCol).Value2
            $ws.Cells.Item($row, $dstCol).Value2 = $srcValue
        }
    }
}
